# Update factsheets with text edits from COMM
#
# The underlying source numbers did not change; only their representation
# did: plain numeric cells became text cells (so Excel stops right-aligning
# / reformatting them), the overall filer count is now shown with a
# thousands separator ("2,849" instead of 2849), and the County sheet grows
# a "Total" summary row (row 89) matching the Total rows already present on
# the other tabs.

function Set-TextValue {
    param($cell, [string]$val)
    # Force Excel to store the value as literal text (no auto number/
    # currency/percent re-interpretation) without leaving a stray
    # NumberFormat/style behind on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overall sheet: A2 (filer count) becomes text with thousands separator
# ---------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")
Set-TextValue $wsOverall.Range("A2") "2,849"

# ---------------------------------------------------------------------
# County sheet: B2:B88 numeric -> text, plus a new Total row (89)
# ---------------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")

$countyCounts = @(
    10,71,18,37,9,6,35,16,11,24,
    22,7,19,22,3,16,6,32,93,2,
    17,8,12,7,21,2,818,14,9,12,
    33,3,6,18,3,10,4,17,3,4,
    4,16,1,6,14,15,5,6,27,16,
    6,14,12,4,71,26,4,12,8,12,
    8,583,2,7,13,42,9,2,17,16,
    10,153,53,18,8,5,9,3,5,5,
    12,53,2,3,25,19,8
)

for ($i = 0; $i -lt $countyCounts.Length; $i++) {
    $row = $i + 2
    $cell = $wsCounty.Cells.Item($row, 2)
    $valStr = [string]$countyCounts[$i]
    Set-TextValue $cell $valStr
}

Set-TextValue $wsCounty.Cells.Item(89, 1) "Total"
Set-TextValue $wsCounty.Cells.Item(89, 2) "2,849"
Set-TextValue $wsCounty.Cells.Item(89, 3) "`$4,691,495,380"
Set-TextValue $wsCounty.Cells.Item(89, 4) "8.08%"
Set-TextValue $wsCounty.Cells.Item(89, 5) "-11.96%"
Set-TextValue $wsCounty.Cells.Item(89, 6) "68.66%"

# ---------------------------------------------------------------------
# Congressional District sheet: B2:B9 numeric -> text; B10 Total -> "2,849"
# ---------------------------------------------------------------------
$wsCd = $wb.Worksheets.Item("Congressional District")

$cdCounts = @(332,141,254,624,605,135,341,417)
for ($i = 0; $i -lt $cdCounts.Length; $i++) {
    $row = $i + 2
    $cell = $wsCd.Cells.Item($row, 2)
    $valStr = [string]$cdCounts[$i]
    Set-TextValue $cell $valStr
}
Set-TextValue $wsCd.Range("B10") "2,849"

# ---------------------------------------------------------------------
# Size sheet: B2:B7 numeric -> text; B8 Total -> "2,849"
# ---------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")

$sizeCounts = @(956,730,406,220,369,168)
for ($i = 0; $i -lt $sizeCounts.Length; $i++) {
    $row = $i + 2
    $cell = $wsSize.Cells.Item($row, 2)
    $valStr = [string]$sizeCounts[$i]
    Set-TextValue $cell $valStr
}
Set-TextValue $wsSize.Range("B8") "2,849"

# ---------------------------------------------------------------------
# Subsector sheet: B2:B13 numeric -> text; B14 Total -> "2,849"
# ---------------------------------------------------------------------
$wsSub = $wb.Worksheets.Item("Subsector")

$subCounts = @(294,335,97,208,54,962,28,11,201,65,566,28)
for ($i = 0; $i -lt $subCounts.Length; $i++) {
    $row = $i + 2
    $cell = $wsSub.Cells.Item($row, 2)
    $valStr = [string]$subCounts[$i]
    Set-TextValue $cell $valStr
}
Set-TextValue $wsSub.Range("B14") "2,849"
